# Update the "metabolomics" omics tag to "Metabolomics" and refresh its
# associated Term Accession Number / Term Source REF so the tag aligns
# with the other (already-NCIT-sourced) tags on the row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

# Tags row: CMML | metabolomics -> Metabolomics | Microbiology | Experiment Metadata
$ws.Range("C13").Value = "Metabolomics"

# Tags Term Accession Number row
$ws.Range("C14").Value = "http://purl.obolibrary.org/obo/NCIT_C49019"

# Tags Term Source REF row
$ws.Range("C15").Value = "NCIT"
